$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new estimate (row 15, column A)
$ws.Range("A15").Value = "agrovoltaics_int_no_involvement_by_population"

# Widen column A so the longer label fits (closest value to 45.43 chars reachable via COM)
$ws.Columns.Item(1).ColumnWidth = 44.666666666666664

# Move the active selection the way the author left it
$ws.Range("A22").Select()
